$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining data row's filename
$ws.Range("A2").Value = "even_MAG-GUT83946.fa"

# Remove rows 3 and 4 (the old row-3/row-4 data is no longer present)
$ws.Range("A3:F4").EntireRow.Delete()
